$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 05:59"

# Pakistan (row 18) - updated case counts
$ws.Range("B18").Value = 154760
$ws.Range("C18").Value = 5839
$ws.Range("D18").Value = 58437
$ws.Range("E18").Value = 93348
$ws.Range("G18").Value = 136
$ws.Range("H18").Value = 2975

# Kazajistan (row 56) - updated case counts
$ws.Range("B56").Value = 15542
$ws.Range("C56").Value = 350
$ws.Range("E56").Value = 5807

# Mongolia (row 163) - updated case counts
$ws.Range("D163").Value = 111
$ws.Range("E163").Value = 86

# Camboya (row 174) - updated case counts
$ws.Range("D174").Value = 126
$ws.Range("E174").Value = 2

# Swap Groenlandia / Islas Malvinas order (rows 206-207)
$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("A207").Value = "Groenlandia"

# Swap Seychelles / Montserrat order (rows 210-211) with updated counts
$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Swap Papua Nueva Guinea / Islas Virgenes Britanicas order (rows 213-214) with updated counts
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
